# "added in some UI Validation"
# Adds a second mini-table below the existing "form field validation" table:
#   row 11 -> blank bold separator cell (A11)
#   row 12 -> section title "UI Validation"
#   row 13 -> column headers
#   rows 14-19 -> one row per social-login UI validation check

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank bold spacer cell (becomes its own cellXfs entry: bold 9pt Helvetica).
$a11 = $ws.Range("A11")
$a11.Font.Bold = $true
$a11.Font.Size = 9
$a11.Font.Name = "Helvetica"

# Section title
$ws.Range("A12").Value = "UI Validation"

# Header row (left to right)
$ws.Range("A13").Value = "Page Name"
$ws.Range("B13").Value = "UI Element"
$ws.Range("C13").Value = "Expected Function"
$ws.Range("D13").Value = "Test ID"
$ws.Range("E13").Value = "Pass/fail"

# Data, filled column-by-column (matches the shared-string insertion order of
# the authored workbook: UI Element list, then Test ID list, then Expected
# Function list, then index/pass columns).
$uiElement = @("Sign in with Facebook", "Sign in with Twitter", "Sign in with Google", "Sign in with Github", "Sign in with Linkedin", "Sign in with Instagram")
$testId    = @("ui.1", "ui.2", "ui.3", "ui.4", "ui.5", "ui.6")
$expected  = @("Successfully login with Facebook", "Successfully login with Twitter", "Successfully login with Google", "Successfully login with Github", "Successfully login with Linkedin", "Successfully login with Instagram ")

for ($i = 0; $i -lt 6; $i++) {
    $ws.Range("B$(14 + $i)").Value = $uiElement[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Range("D$(14 + $i)").Value = $testId[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Range("C$(14 + $i)").Value = $expected[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Range("A$(14 + $i)").Value = "index"
    $ws.Range("E$(14 + $i)").Value = "pass"
}

# Column widths for the new columns (B gets created, C gets widened)
$ws.Columns(2).ColumnWidth = 18.3
$ws.Columns(3).ColumnWidth = 26.6

# Match the final selection left behind in the saved file
$ws.Range("E20").Select() | Out-Null
